$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# are first marked as Text (@) so the literal string survives round-tripping.

$ws.Range("D2").Value2 = '30.212.35'
$ws.Range("E2").Value2 = '  +0.61%  '
$ws.Range("D3").Value2 = '1.912.93'
$ws.Range("E3").Value2 = '  +0.12%  '
$ws.Range("E4").Value2 = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '0.8178'
$ws.Range("E5").Value2 = '  +3.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '243.60'
$ws.Range("E6").Value2 = '  +0.37%  '
$ws.Range("E7").Value2 = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '0.3253'
$ws.Range("E8").Value2 = '  +2.92%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '26.79'
$ws.Range("E9").Value2 = '  +2.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '0.07065'
$ws.Range("E10").Value2 = '  +2.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '0.08093'
$ws.Range("E11").Value2 = '  +1.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '0.7794'
$ws.Range("E12").Value2 = '  +4.37%  '
$ws.Range("D13").Value2 = '1.900.08'
$ws.Range("E13").Value2 = '  -0.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '5.302'
$ws.Range("E14").Value2 = '  +1.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '93.38'
$ws.Range("E15").Value2 = '  +0.17%  '
$ws.Range("D16").Value2 = '30.207.48'
$ws.Range("E16").Value2 = '  +0.51%  '
$ws.Range("E17").Value2 = '  +1.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '5.927'
$ws.Range("E18").Value2 = '  -0.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '247.76'
$ws.Range("E19").Value2 = '  +0.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '0.000007805'
$ws.Range("E20").Value2 = '  +0.29%  '
$ws.Range("D21").Value2 = '2.168.05'
$ws.Range("E21").Value2 = '  +0.38%  '
$ws.Range("E22").Value2 = '  +0.06%  '
$ws.Range("E23").Value2 = '  -0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '7.248'
$ws.Range("E24").Value2 = '  +5.08%  '
$ws.Range("E25").Value2 = '  +21.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '9.350'
$ws.Range("E26").Value2 = '  +0.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '167.67'
$ws.Range("E27").Value2 = '  -1.23%  '
$ws.Range("E28").Value2 = '  +0.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '2.118'
$ws.Range("E29").Value2 = '  +3.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '1.373'
$ws.Range("E30").Value2 = '  -0.51%  '
$ws.Range("E31").Value2 = '  +0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '4.317'
$ws.Range("E32").Value2 = '  -0.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '0.05820'
$ws.Range("E33").Value2 = '  +5.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '4.108'
$ws.Range("E34").Value2 = '  -0.01%  '
$ws.Range("E35").Value2 = '  +1.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '0.7382'
$ws.Range("E36").Value2 = '  +0.25%  '
$ws.Range("B37").Value2 = 'Frax'
$ws.Range("C37").Value2 = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '1.003'
$ws.Range("E37").Value2 = '  +0.37%  '
$ws.Range("B38").Value2 = 'HuobiToken'
$ws.Range("C38").Value2 = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '2.706'
$ws.Range("E38").Value2 = '  -0.79%  '
$ws.Range("B39").Value2 = 'VeChain'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '0.01930'
$ws.Range("E39").Value2 = '  -0.73%  '
$ws.Range("B40").Value2 = 'MXToken'
$ws.Range("C40").Value2 = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '2.800'
$ws.Range("E40").Value2 = '  +0.16%  '
$ws.Range("B41").Value2 = 'TheSandbox'
$ws.Range("C41").Value2 = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '0.4471'
$ws.Range("E41").Value2 = '  +0.79%  '
$ws.Range("B42").Value2 = 'Aave'
$ws.Range("C42").Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '73.46'
$ws.Range("E42").Value2 = '  +1.38%  '
$ws.Range("B43").Value2 = 'FraxShare'
$ws.Range("C43").Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '5.968'
$ws.Range("E43").Value2 = '  -3.42%  '
$ws.Range("B44").Value2 = 'TrustWalletToken'
$ws.Range("C44").Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '0.8519'
$ws.Range("E44").Value2 = '  +1.93%  '
$ws.Range("B45").Value2 = 'RenderToken'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '1.914'
$ws.Range("E45").Value2 = '  +0.92%  '
$ws.Range("B46").Value2 = 'PaxDollar'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '1.001'
$ws.Range("E46").Value2 = '  -0.06%  '
$ws.Range("B47").Value2 = 'Quant'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '102.97'
$ws.Range("E47").Value2 = '  +2.44%  '
$ws.Range("B48").Value2 = 'Maker'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value2 = '1.015.51'
$ws.Range("E48").Value2 = '  +3.14%  '
$ws.Range("B49").Value2 = 'Aptos'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '7.600'
$ws.Range("E49").Value2 = '  +0.63%  '
$ws.Range("B50").Value2 = 'EnergySwap'
$ws.Range("C50").Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '9.879'
$ws.Range("E50").Value2 = '  +0.93%  '
$ws.Range("B51").Value2 = 'RocketPoolETH'
$ws.Range("C51").Value2 = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value2 = '2.062.55'
$ws.Range("E51").Value2 = '  +0.13%  '
